$wb = $excel.ActiveWorkbook

# ---------- Sheet "VENTAS POR GRUPO" (columns A:R) ----------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Duplicate row 16 (SOLIS OCAMPO DIMAS ABDON, with its 43.86 LAVABOS sale) and
# insert the copy as the new row 17, preserving formatting. This pushes the
# old rows 17 (VEHINVER SA) and 18 (totals) down to 18 and 19.
$ws1.Rows.Item(16).Copy()
$ws1.Rows.Item(17).Insert()

# Row 16 becomes a new asesor "QUIJIJE MENDOZA GENESIS XIOMARA" with no sales yet
$ws1.Range("B16").Value = "QUIJIJE MENDOZA GENESIS XIOMARA"
$ws1.Range("I16").Value = 0

# Update the totals row (now row 19) counters from "de 16" to "de 17"
$ws1.Range("C19").Value = "0 de 17"
$ws1.Range("D19").Value = "0 de 17"
$ws1.Range("E19").Value = "1 de 17"
$ws1.Range("F19").Value = "0 de 17"
$ws1.Range("G19").Value = "0 de 17"
$ws1.Range("H19").Value = "0 de 17"
$ws1.Range("I19").Value = "1 de 17"
$ws1.Range("J19").Value = "0 de 17"
$ws1.Range("K19").Value = "0 de 17"
$ws1.Range("L19").Value = "3 de 17"
$ws1.Range("M19").Value = "1 de 17"
$ws1.Range("N19").Value = "0 de 17"
$ws1.Range("O19").Value = "1 de 17"
$ws1.Range("P19").Value = "0 de 17"
$ws1.Range("Q19").Value = "0 de 17"
$ws1.Range("R19").Value = "0 de 17"

# ---------- Sheet "VENTA MENSUAL" (columns A:G) ----------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Duplicate row 16 (SOLIS OCAMPO DIMAS ABDON, with its 43.86 julio sale) and
# insert the copy as the new row 17, preserving formatting. This pushes the
# old rows 17 (VEHINVER SA) and 18 (totals) down to 18 and 19.
$ws2.Rows.Item(16).Copy()
$ws2.Rows.Item(17).Insert()

# Row 16 becomes a new asesor "QUIJIJE MENDOZA GENESIS XIOMARA" with no sales yet
$ws2.Range("B16").Value = "QUIJIJE MENDOZA GENESIS XIOMARA"
$ws2.Range("F16").Value = 0
